# Add Files via upload
# Appends the two newest HackerRank "Algorithms" progress-log entries to the
# "Math & Algorithms" sheet's second table (Table4, columns I:N), extends
# the table range to include them, and leaves the same cell selected that
# the author had selected when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Math & Algorithms")

# --- Copy the formatting of the last existing table row (19) down onto the
#     two new rows (20-21) so the new rows look like part of the table. ---
$ws.Range("I19:N19").Copy()
$ws.Range("I20:N21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 20 : 2023-07-27 -----------------------------------------------
$ws.Range("I20").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("J20").Value = 45134
$ws.Range("L20").Value = "1712.97/2200"
$ws.Range("K20").Value = "64% (487.03 points to next star)"
$ws.Range("M20").Value = 66794
$ws.Range("N20").Formula = "=IF(ROW()>2,(`$M`$2-M20)/`$M`$2,`"NA`")"

# --- Row 21 : 2023-07-28 -----------------------------------------------
$ws.Range("I21").Value = "Problem Solving(Algorithms & Data Structures)"
$ws.Range("J21").Value = 45135
$ws.Range("L21").Value = "1732.97/2200"
$ws.Range("K21").Value = "65% (467.03 points to next star)"
$ws.Range("M21").Value = 65489
$ws.Range("N21").Formula = "=IF(ROW()>2,(`$M`$2-M21)/`$M`$2,`"NA`")"

# --- Grow "Table4" (the I1:N19 structured table) so it covers the two
#     freshly-added rows, mirroring what Excel does when a table is
#     extended. ---
$lo = $ws.ListObjects.Item("Table4")
$lo.Resize($ws.Range("I1:N21"))

# --- Match the author's final selection/view state. ---
[void]$ws.Range("I20:I21").Select()
